$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" '68.057.08'
$ws.Range("E2").Value = '  +1.78%  '
Set-TextValue "D3" '3.907.79'
$ws.Range("E4").Value = '  +0.21%  '
Set-TextValue "D5" '484.58'
$ws.Range("E5").Value = '  +3.67%  '
Set-TextValue "D6" '145.68'
$ws.Range("E6").Value = '  +0.28%  '
Set-TextValue "D7" '0.624'
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("E8").Value = '  -0.15%  '
Set-TextValue "D9" '0.727'
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("E10").Value = '  +1.54%  '
Set-TextValue "D11" '0.0000353'
$ws.Range("E11").Value = '  +4.67%  '
Set-TextValue "D12" '42.41'
$ws.Range("E12").Value = '  -1.94%  '
Set-TextValue "D13" '10.60'
$ws.Range("E13").Value = '  +1.47%  '
Set-TextValue "D14" '4.531.87'
$ws.Range("E14").Value = '  +0.08%  '
Set-TextValue "D15" '14.69'
$ws.Range("E15").Value = '  -2.29%  '
Set-TextValue "D16" '3.908.79'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("E17").Value = '  -0.13%  '
Set-TextValue "D18" '19.78'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("E19").Value = '  -2.48%  '
Set-TextValue "D20" '68.138.74'
$ws.Range("E20").Value = '  +1.44%  '
Set-TextValue "D21" '447.94'
$ws.Range("E21").Value = '  +3.82%  '
Set-TextValue "D22" '14.68'
$ws.Range("E22").Value = '  -0.33%  '
Set-TextValue "D23" '3.34'
$ws.Range("E23").Value = '  -0.39%  '
Set-TextValue "D24" '88.91'
$ws.Range("E24").Value = '  +0.41%  '
Set-TextValue "D25" '11.58'
$ws.Range("E25").Value = '  +14.75%  '
Set-TextValue "D26" '11.10'
$ws.Range("E26").Value = '  +14.05%  '
Set-TextValue "D27" '3.60'
$ws.Range("E27").Value = '  +2.17%  '
Set-TextValue "D28" '38.72'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D30" '690.06'
$ws.Range("E30").Value = '  -6.55%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D31" '13.36'
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D32" '0.130'
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("E33").Value = '  +2.84%  '
Set-TextValue "D34" '0.0₃0923'
$ws.Range("E34").Value = '  +22.33%  '
Set-TextValue "D35" '41.62'
$ws.Range("E35").Value = '  -5.16%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D36" '58.89'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D37" '5.72'
$ws.Range("E37").Value = '  +6.56%  '
Set-TextValue "D38" '0.149'
$ws.Range("E38").Value = '  -5.54%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D40" '2.91'
$ws.Range("E40").Value = '  +16.94%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D41" '0.0477'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D42" '0.364'
$ws.Range("E42").Value = '  +8.01%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue "D43" '3.04'
$ws.Range("E43").Value = '  -6.74%  '
Set-TextValue "D44" '3.01'
$ws.Range("E44").Value = '  +7.24%  '
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("E46").Value = '  -0.04%  '
Set-TextValue "D47" '3.41'
$ws.Range("E47").Value = '  -0.93%  '
Set-TextValue "D48" '2.12'
$ws.Range("E48").Value = '  -3.11%  '
Set-TextValue "D49" '146.45'
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D51" '2.82'
$ws.Range("E51").Value = '  -1.90%  '
